# Update the per-row runs/balls/fours/sixes stats (columns C:F) on the
# "Ishan Kishan " sheet to reflect the reshuffled innings figures.
# Row 4 (72/47/8/3) is unchanged. Values are written as text ("@" number
# format) to preserve the original numberStoredAsText string typing.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @("33","19","3","1")
    3  = @("55","30","4","3")
    5  = @("7","7","1","0")
    6  = @("68","37","6","5")
    7  = @("99","58","2","9")
    8  = @("33","30","1","2")
    9  = @("28","15","2","2")
    10 = @("37","36","4","1")
    11 = @("28","32","1","1")
    12 = @("0","1","0","0")
    13 = @("25","19","3","1")
    14 = @("31","23","1","2")
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt 4; $i++) {
        $col = 3 + $i
        $cell = $ws.Cells.Item($row, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $vals[$i]
    }
}
